# Christos completed his skills
# Fill in Christos's column (D) on the SkillsMatrix sheet with his
# self-assessment ratings, using the same S#/E#/SE#/"-" vocabulary as
# the other team members' columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SkillsMatrix")

# [Programming Languages]
$ws.Range("D3").Value  = "S2"   # Java SE
$ws.Range("D4").Value  = "-"    # Java EE
$ws.Range("D5").Value  = "-"    # Java ME
$ws.Range("D6").Value  = "-"    # Android
$ws.Range("D7").Value  = "S2"   # C++
$ws.Range("D8").Value  = "E3"   # C#
$ws.Range("D9").Value  = "-"    # JavaScript
$ws.Range("D10").Value = "S1"   # HTML
$ws.Range("D11").Value = "S3"   # SQL

# [Development Processes]
$ws.Range("D13").Value = "-"    # Agile
$ws.Range("D14").Value = "-"    # Waterfall

# [Development Tools]
$ws.Range("D16").Value = "S1"   # JDeveloper (UML)
$ws.Range("D17").Value = "S3"   # NetBeans
$ws.Range("D18").Value = "S1"   # Eclipse
$ws.Range("D19").Value = "-"    # Ant
$ws.Range("D20").Value = "-"    # Maven
$ws.Range("D21").Value = "-"    # GitHub
$ws.Range("D22").Value = "-"    # Amazon ES2
$ws.Range("D23").Value = "-"    # MS Visual Source

# [Project Management]
$ws.Range("D25").Value = "S2"   # Project
$ws.Range("D26").Value = "S1"   # MS Project (Gantt charts)
$ws.Range("D27").Value = "E3"   # Microsoft Office

# [Product Development]
$ws.Range("D29").Value = "SE3"  # User Interaction Design
$ws.Range("D30").Value = "S2"   # Research Methods

# [Languages]
$ws.Range("D34").Value = "S3"   # English

# Leave the view the way it ended up after this editing pass.
$ws.Range("D10").Select()
